$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row: A2 = "C123413", B2 = 2
$ws.Range("A2").Value = "C123413"
$ws.Range("B2").Value = 2

# Match the author's final selection state (row 3 selected)
$ws.Range("A3:XFD3").Select()
